$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 45: 2025-11-19
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2025-11-19"
$ws.Range("A45").ClearFormats()
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = 26

# Row 46: 2025-11-20
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "2025-11-20"
$ws.Range("A46").ClearFormats()
$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 25
